# Add two new worksheets ("MDN" and "Google") with resource links, at the
# end of the workbook, and leave "Google" as the active/selected sheet
# (mirrors the author's commit: "Google and MDN Resourses added.")

$wb = $excel.ActiveWorkbook

# --- Add "MDN" worksheet after the last existing sheet (Front-End Master) ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsMdn = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsMdn.Name = "MDN"

$wsMdn.Range("C2").Value = "https://developer.mozilla.org/en-US/docs/Web/HTML"

$wsMdn.Range("C4").Value = "https://developer.mozilla.org/en-US/docs/Web/HTML/Element/button"
$wsMdn.Range("B4").Value = "Buttons"

$wsMdn.Range("C6").Value = "https://developer.mozilla.org/en-US/docs/Web/HTML/Element/figure"
$wsMdn.Range("B6").Value = "Figure"

$wsMdn.Range("C8").Value = "https://developer.mozilla.org/en-US/docs/Web/HTML/Element"
$wsMdn.Range("B8").Value = "HTML Reference"

$wsMdn.Range("C10").Value = "https://developer.mozilla.org/en-US/docs/Web/CSS/Pseudo-classes"
$wsMdn.Range("B10").Value = "Pseudo Classes CSS"

# NOTE: the host's ColumnWidth setter quantizes to whole-pixel increments,
# so the literal target widths (89.109375 / 88.44140625) can't be hit
# exactly. These inputs land on the closest achievable stored width.
$wsMdn.Columns.Item(2).ColumnWidth = 88.252
$wsMdn.Columns.Item(3).ColumnWidth = 87.586

$wsMdn.Range("C12").Select() | Out-Null

# --- Add "Google" worksheet after "MDN" (i.e. at the very end) ---
$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsGoogle = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet2)
$wsGoogle.Name = "Google"

$wsGoogle.Range("C2").Value = "https://www.chromestatus.com/samples"
$wsGoogle.Range("B2").Value = "Chrome Samples"

$wsGoogle.Range("C4").Value = "https://developers.google.com/web/fundamentals"
$wsGoogle.Range("B4").Value = "Web Fundamentals"

# Same pixel-quantization caveat as above.
$wsGoogle.Columns.Item(2).ColumnWidth = 61.252
$wsGoogle.Columns.Item(3).ColumnWidth = 88.15

$wsGoogle.Range("B6").Select() | Out-Null

# "Google" ends up as the active tab, matching activeTab="4" in workbook.xml
$wsGoogle.Activate()
